$d = $word.ActiveDocument

$replacements = @(
    @{old = "108×8="; new = "262×5="},
    @{old = "516×6="; new = "226×3="},
    @{old = "897×4="; new = "111×6="},
    @{old = "528×6="; new = "187×6="},
    @{old = "919×4="; new = "903×6="},
    @{old = "858×2="; new = "300×3="},
    @{old = "534×9="; new = "789×6="},
    @{old = "692×3="; new = "896×2="},
    @{old = "993×4="; new = "990×5="},
    @{old = "885×4="; new = "934×8="},
    @{old = "779×7="; new = "676×3="},
    @{old = "725×3="; new = "389×3="},
    @{old = "688×3="; new = "992×7="},
    @{old = "279×5="; new = "584×3="},
    @{old = "927×9="; new = "194×5="},
    @{old = "564×3="; new = "647×9="},
    @{old = "800×5="; new = "955×2="},
    @{old = "543×6="; new = "489×8="},
    @{old = "101×6="; new = "482×5="},
    @{old = "979×5="; new = "528×4="},
    @{old = "385×5="; new = "483×2="},
    @{old = "426×9="; new = "258×4="},
    @{old = "198×5="; new = "600×4="},
    @{old = "284×4="; new = "822×6="},
    @{old = "558×7="; new = "128×2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
